$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Cells changing type/style (copy formatting+value from a donor cell with identical target content) ---
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("D14").Copy($ws.Range("D20"))
$ws.Range("H14").Copy($ws.Range("E20"))
$ws.Range("F14").Copy($ws.Range("D22"))
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("G14").Copy($ws.Range("C23"))
$ws.Range("C14").Copy($ws.Range("G26"))
$ws.Range("L14").Copy($ws.Range("H26"))
$ws.Range("D14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("N14").Copy($ws.Range("M28"))
$ws.Range("N15").Copy($ws.Range("M29"))
$ws.Range("F14").Copy($ws.Range("D30"))
$ws.Range("H14").Copy($ws.Range("E30"))
$ws.Range("K27").Copy($ws.Range("L30"))

# --- Simple numeric value changes (style unchanged) ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = -57.142857142857
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = -78.571428571428
$ws.Range("N16").Value = -96
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 5
$ws.Range("J17").Value = 6
$ws.Range("K17").Value = -16.666666666666
$ws.Range("L17").Value = -54.545454545454
$ws.Range("M17").Value = -28.571428571428
$ws.Range("N17").Value = -82.142857142857
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 18
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = -73.333333333333
$ws.Range("L18").Value = -42.857142857142
$ws.Range("M18").Value = -42.857142857142
$ws.Range("N18").Value = -91.836734693877
$ws.Range("C19").Value = 2
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = -6.25
$ws.Range("I19").Value = 10
$ws.Range("J19").Value = 11
$ws.Range("K19").Value = -9.090909090909
$ws.Range("L19").Value = -9.090909090909
$ws.Range("M19").Value = -23.076923076923
$ws.Range("N19").Value = -61.538461538461
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 13
$ws.Range("K20").Value = 225
$ws.Range("L20").Value = 44.444444444444
$ws.Range("M20").Value = 225
$ws.Range("N20").Value = -72.340425531914
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = -38.461538461538
$ws.Range("F21").Value = 52
$ws.Range("G21").Value = 54
$ws.Range("H21").Value = -3.703703703703
$ws.Range("I21").Value = 35
$ws.Range("J21").Value = 43
$ws.Range("K21").Value = -18.604651162790
$ws.Range("L21").Value = -16.666666666666
$ws.Range("M21").Value = -22.222222222222
$ws.Range("N21").Value = -84.782608695652
$ws.Range("M22").Value = -50
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 66.666666666666
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -60
$ws.Range("M23").Value = -77.777777777777
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 43
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = -12.244897959183
$ws.Range("I24").Value = 32
$ws.Range("J24").Value = 37
$ws.Range("K24").Value = -13.513513513513
$ws.Range("L24").Value = -27.272727272727
$ws.Range("M24").Value = -33.333333333333
$ws.Range("C25").Value = 8
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 83.333333333333
$ws.Range("I25").Value = 26
$ws.Range("J25").Value = 13
$ws.Range("K25").Value = 100
$ws.Range("L25").Value = 18.181818181818
$ws.Range("M25").Value = 136.363636363636
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0
